# Fix weekly forecast data: shift all week-start dates forward by one week
# (the previous save failed to roll the forecast window forward) and
# refresh the Amazon P90 forecast values (and one P80 value) that changed
# as a result, plus the dependent Summary-sheet figures.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date values for rows 2-17 (one week later than before).
# Force text format first so Excel does not silently convert these
# yyyy-mm-dd strings into date serial numbers.
$dateRange = $wsForecast.Range("B2:B17")
$dateRange.NumberFormat = "@"

$newDates = @(
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27",
    "2025-05-04",
    "2025-05-11",
    "2025-05-18"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Range("B$row").Value = $newDates[$i]
}

# Updated Amazon P90 Forecast values (column H) for the affected rows.
$wsForecast.Range("H2").Value  = 4
$wsForecast.Range("H8").Value  = 4
$wsForecast.Range("H12").Value = 3
$wsForecast.Range("H13").Value = 3
$wsForecast.Range("H15").Value = 3

# Row 17 also has an updated Amazon P80 Forecast (column G) and a
# correspondingly updated Amazon P90 Forecast (column H).
$wsForecast.Range("G17").Value = 1
$wsForecast.Range("H17").Value = 2

# Summary sheet: historical range now extends one week further, and the
# max/min forecast week references move forward one week as well.
# Force text format on the week-reference cells so the yyyy-mm-dd strings
# are not silently converted into date serial numbers.
$wsSummary.Range("B13:B15").NumberFormat = "@"

$wsSummary.Range("B2").Value  = "2022-12-25 to 2025-01-26"
$wsSummary.Range("B13").Value = "2025-02-02"
$wsSummary.Range("B15").Value = "2025-02-02"
